# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" detail table (rows 16-56, columns E:G) is refreshed:
#   - Column E ("Periodo Mora") is re-sorted ascending (1610 .. 2003) instead
#     of the previous descending order (2003 .. 1610).
#   - Column F ("Valor Mora") keeps its two-tier value but the tiers swap
#     which periods they apply to: the first 22 periods (1610-1808) now get
#     27578 and the remaining 19 periods (1809-2003) get 31249.
#   - Column G ("Salario Basico") is updated uniformly from 689455 to 781242.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Ascending list of the 41 "Periodo Mora" labels that now populate E16:E56.
$periodos = @(
    "1610","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$firstDataRow = 16
$newSalarioBasico = 781242

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = $firstDataRow + $i
    $periodo = $periodos[$i]

    # Periods through 1808 (the first 22 rows) move to the lower "Valor Mora"
    # tier; the rest (1809 onward) move to the higher tier.
    if ($i -lt 22) {
        $valorMora = 27578
    } else {
        $valorMora = 31249
    }

    $ws.Range("E$row").Value = $periodo
    $ws.Range("F$row").Value = $valorMora
    $ws.Range("G$row").Value = $newSalarioBasico
}
